$wb = $excel.ActiveWorkbook

# --- Add the new "Com" worksheet at the end of the workbook ---
$wsCom = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCom.Name = "Com"

# Header row (set B1 before A1 so new shared strings are appended in the
# same order as the target workbook: BldgLoc=34, BldgType=35)
$wsCom.Range("B1").Value = "BldgLoc"
$wsCom.Range("A1").Value = "BldgType"
$wsCom.Range("C1").Value = "Normunit"
$wsCom.Range("D1").Value = "Value"
$wsCom.Range("E1").Value = "BldgVint"

# Data rows: BldgType code, Normunit (Area-ft2), Value
$rows = @(
  @("Asm", 100002.1),
  @("ECC", 299999.59999999998),
  @("EPr", 50000.14),
  @("ERC", 1920.0170000000001),
  @("ESe", 149998.6),
  @("EUn", 930201.4),
  @("Hsp", 235501),
  @("Htl", 139998.9),
  @("MBT", 199999.2),
  @("MLI", 100001.1),
  @("Mtl", 30000.07),
  @("Nrs", 60654.58),
  @("OfL", 174998.9),
  @("OfS", 10000.08),
  @("RFF", 2500.0529999999999),
  @("RSD", 5599.9570000000003),
  @("Rt3", 120000.5),
  @("RtL", 129997),
  @("RtS", 7999.9290000000001),
  @("SCn", 250000.3)
)

$r = 2
foreach ($row in $rows) {
    $wsCom.Cells.Item($r, 1).Value = $row[0]
    $wsCom.Cells.Item($r, 3).Value = "Area-ft2-BA"
    $wsCom.Cells.Item($r, 4).Value = $row[1]
    $r = $r + 1
}

$wsCom.Columns.Item(3).ColumnWidth = 11.75

# --- Adjust selections on the other sheets to match the saved view state ---
$wsDMo = $wb.Worksheets.Item("DMo")
[void]$wsDMo.Range("A2").Select()

$wsSFm = $wb.Worksheets.Item("SFm")
[void]$wsSFm.Range("D50").Select()

# --- Com tab is the active / selected tab when the file was saved ---
[void]$wsCom.Select()
[void]$wsCom.Range("G6").Select()
